$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B66:B82 Testable column from "y" to "n" ---
$ws.Range("B66:B82").Value = 'n'

# --- Append new rows 83-90 for map/multiset test cases ---
# Row 83
$ws.Range('A83').Value = 'updel_082'
$ws.Range('B83').Value = 'y'
$ws.Range('C83').Value = '删除map字段值为null的数据可删除成功'
$ws.Range('D83').Value = 'ComplexDataType'
$ws.Range('E83').Value = 'Map'
$ws.Range('F83').Value = 'map1'
$ws.Range('G83').Value = 'map1_value07'
$ws.Range('H83').Value = 'delete from $map1 where id=10'
$ws.Range('I83').Value = '1'
$ws.Range('J83').Value = 'select * from $map1'
$ws.Range('K83').Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_082.csv'
$ws.Range('L83').Value = 'csv_equals'

# Row 84
$ws.Range('A84').Value = 'updel_083'
$ws.Range('B84').Value = 'y'
$ws.Range('C84').Value = '删除含有多个map字段的表数据'
$ws.Range('D84').Value = 'ComplexDataType'
$ws.Range('E84').Value = 'Map'
$ws.Range('F84').Value = 'map2'
$ws.Range('G84').Value = 'map2_value13'
$ws.Range('H84').Value = 'delete from $map2 where id=2'
$ws.Range('I84').Value = '1'
$ws.Range('J84').Value = 'select * from $map2'
$ws.Range('K84').Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_083.csv'
$ws.Range('L84').Value = 'csv_containsAll'

# Row 85
$ws.Range('A85').Value = 'updel_084'
$ws.Range('B85').Value = 'y'
$ws.Range('C85').Value = '表中含有map类型字段，更新整型字段值'
$ws.Range('D85').Value = 'ComplexDataType'
$ws.Range('E85').Value = 'Map'
$ws.Range('F85').Value = 'map7'
$ws.Range('G85').Value = 'map7_value20'
$ws.Range('H85').Value = 'update $map7 set age=20 where id=1'
$ws.Range('I85').Value = '1'
$ws.Range('J85').Value = 'select * from $map7'
$ws.Range('K85').Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_084.csv'
$ws.Range('L85').Value = 'csv_containsAll'

# Row 86
$ws.Range('A86').Value = 'updel_085'
$ws.Range('B86').Value = 'y'
$ws.Range('C86').Value = '表中含有map类型字段，更新字符型字段值'
$ws.Range('D86').Value = 'ComplexDataType'
$ws.Range('E86').Value = 'Map'
$ws.Range('F86').Value = 'map7'
$ws.Range('G86').Value = 'map7_value20'
$ws.Range('H86').Value = 'update $map7 set name=''dingo'''
$ws.Range('I86').Value = '2'
$ws.Range('J86').Value = 'select id,name,data from $map7'
$ws.Range('K86').Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_085.csv'
$ws.Range('L86').Value = 'csv_containsAll'

# Row 87
$ws.Range('A87').Value = 'updel_086'
$ws.Range('B87').Value = 'n'
$ws.Range('C87').Value = '删除含有multiset类型字段的单条数据'
$ws.Range('D87').Value = 'ComplexDataType'
$ws.Range('E87').Value = 'Multiset'
$ws.Range('F87').Value = 'multiset20'
$ws.Range('G87').Value = 'multiset20_value47'
$ws.Range('H87').Value = 'delete from $multiset20 where id=3'
$ws.Range('I87').Value = '1'
$ws.Range('J87').Value = 'select * from $multiset20'
$ws.Range('K87').Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_086.csv'
$ws.Range('L87').Value = 'csv_containsAll'

# Row 88
$ws.Range('A88').Value = 'updel_087'
$ws.Range('B88').Value = 'n'
$ws.Range('C88').Value = '删除含有multiset类型字段的全表数据'
$ws.Range('D88').Value = 'ComplexDataType'
$ws.Range('E88').Value = 'Multiset'
$ws.Range('F88').Value = 'multiset20'
$ws.Range('G88').Value = 'multiset20_value47'
$ws.Range('H88').Value = 'delete from $multiset20'
$ws.Range('I88').Value = '3'
$ws.Range('J88').Value = 'select * from $multiset20'
$ws.Range('K88').Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_087.csv'
$ws.Range('L88').Value = 'csv_equals'

# Row 89
$ws.Range('A89').Value = 'updel_088'
$ws.Range('B89').Value = 'n'
$ws.Range('C89').Value = '删除整型列值为null的行数据'
$ws.Range('D89').Value = 'ComplexDataType'
$ws.Range('E89').Value = 'Multiset'
$ws.Range('F89').Value = 'multiset1'
$ws.Range('G89').Value = 'multiset1_value20'
$ws.Range('H89').Value = 'delete from $multiset1'
$ws.Range('I89').Value = '1'
$ws.Range('J89').Value = 'select * from $multiset1'
$ws.Range('K89').Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_088.csv'
$ws.Range('L89').Value = 'csv_equals'

# Row 90
$ws.Range('A90').Value = 'updel_089'
$ws.Range('B90').Value = 'n'
$ws.Range('C90').Value = '删除字符型列值为null的行数据'
$ws.Range('D90').Value = 'ComplexDataType'
$ws.Range('E90').Value = 'Multiset'
$ws.Range('F90').Value = 'multiset3'
$ws.Range('G90').Value = 'multiset3_value21'
$ws.Range('H90').Value = 'delete from $multiset3 where id=1'
$ws.Range('I90').Value = '1'
$ws.Range('J90').Value = 'select * from $multiset3'
$ws.Range('K90').Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_089.csv'
$ws.Range('L90').Value = 'csv_equals'

# --- Update the sheet view selection (scrolled viewport itself is not
#     part of persisted model state in this runtime, only the selection is) ---
$ws.Range('B82').Select() | Out-Null
